$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2200.1667
$ws.Range("J70").Value = 3253
$ws.Range("L70").Value = 9759
$ws.Range("N70").Value = -10299
$ws.Range("H73").Value = 2200.1667
$ws.Range("J73").Value = 3253
$ws.Range("L73").Value = 9759
$ws.Range("N73").Value = -11631
$ws.Range("H98").Value = 2664.8215
$ws.Range("I98").Value = 2664.8215
$ws.Range("K98").Value = 2664.8215
$ws.Range("M98").Value = -1166.8215
$ws.Range("H122").Value = 2664.8215
$ws.Range("I122").Value = 2664.8215
$ws.Range("K122").Value = 7994.4645
$ws.Range("M122").Value = -5544.4645
$ws.Range("H132").Value = 12827575
$ws.Range("I132").Value = 15159217
$ws.Range("J132").Value = 3544
$ws.Range("K132").Value = 45477651
$ws.Range("L132").Value = 10632
$ws.Range("M132").Value = -45475121
$ws.Range("N132").Value = -15692
$ws.Range("H135").Value = 1082.64
$ws.Range("I135").Value = 283.73685
$ws.Range("K135").Value = 2553.63165
$ws.Range("M135").Value = -18.63165000000026

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1578.1666
$ws.Range("I45").Value = 2019.7142
$ws.Range("J45").Value = 960
$ws.Range("K45").Value = 2019.7142
$ws.Range("L45").Value = 960
$ws.Range("M45").Value = -1642.7142
$ws.Range("N45").Value = -1714
$ws.Range("H61").Value = 799.25
$ws.Range("I61").Value = 799.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 799.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -587.25
$ws.Range("N61").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H124").Value = 3500
$ws.Range("J124").Value = 3500
$ws.Range("L124").Value = 3500
$ws.Range("N124").Value = -13320
$ws.Range("H125").Value = 38000
$ws.Range("J125").Value = 38000
$ws.Range("L125").Value = 38000
$ws.Range("N125").Value = -47840
$ws.Range("H132").Value = 1752.7188
$ws.Range("I132").Value = 1442.826
$ws.Range("K132").Value = 4328.478
$ws.Range("M132").Value = -1798.478
$ws.Range("H133").Value = 28618.234
$ws.Range("J133").Value = 28618.234
$ws.Range("L133").Value = 28618.234
$ws.Range("N133").Value = -33678.234
$ws.Range("H135").Value = 20500
$ws.Range("J135").Value = 20500
$ws.Range("L135").Value = 20500
$ws.Range("N135").Value = -30640
$ws.Range("H136").Value = 799.25
$ws.Range("I136").Value = 799.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2397.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 152.25
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 34790
$ws.Range("J139").Value = 34790
$ws.Range("L139").Value = 34790
$ws.Range("N139").Value = -45070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6453.9565
$ws.Range("I134").Value = 970.6316
$ws.Range("K134").Value = 2911.8948
$ws.Range("M134").Value = -376.8948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1413.5834
$ws.Range("I31").Value = 1133.9412
$ws.Range("K31").Value = 1133.9412
$ws.Range("M31").Value = -838.9412
$ws.Range("H34").Value = 1413.5834
$ws.Range("I34").Value = 1133.9412
$ws.Range("K34").Value = 1133.9412
$ws.Range("M34").Value = -931.9412
$ws.Range("H99").Value = 2084.75
$ws.Range("J99").Value = 2399
$ws.Range("L99").Value = 2399
$ws.Range("N99").Value = -5395
$ws.Range("H107").Value = 578.35
$ws.Range("I107").Value = 444.2143
$ws.Range("K107").Value = 444.2143
$ws.Range("M107").Value = 1475.7857
$ws.Range("H122").Value = 1184.6666
$ws.Range("I122").Value = 1018
$ws.Range("J122").Value = 1518
$ws.Range("K122").Value = 3054
$ws.Range("L122").Value = 4554
$ws.Range("M122").Value = -604
$ws.Range("N122").Value = -9454
$ws.Range("H126").Value = 2084.75
$ws.Range("J126").Value = 2399
$ws.Range("L126").Value = 7197
$ws.Range("N126").Value = -12137
$ws.Range("H132").Value = 7138.25
$ws.Range("I132").Value = 9547.416999999999
$ws.Range("K132").Value = 28642.251
$ws.Range("M132").Value = -26112.251

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 38.416668
$ws.Range("I2").Value = 12.666667
$ws.Range("J2").Value = 115.666664
$ws.Range("K2").Value = 76.00000199999999
$ws.Range("L2").Value = 693.999984
$ws.Range("M2").Value = 36.99999800000001
$ws.Range("N2").Value = -919.999984
$ws.Range("H80").Value = 4498.625
$ws.Range("J80").Value = 4498.625
$ws.Range("L80").Value = 13495.875
$ws.Range("N80").Value = -15367.875
$ws.Range("H83").Value = 4498.625
$ws.Range("J83").Value = 4498.625
$ws.Range("L83").Value = 40487.625
$ws.Range("N83").Value = -49847.625
$ws.Range("H131").Value = 10753811
$ws.Range("J131").Value = 1151.5222
$ws.Range("L131").Value = 3454.5666
$ws.Range("N131").Value = -13534.5666
$ws.Range("H136").Value = 1327.2
$ws.Range("I136").Value = 1271.1111
$ws.Range("J136").Value = 1832
$ws.Range("K136").Value = 3813.3333
$ws.Range("L136").Value = 5496
$ws.Range("M136").Value = 1286.6667
$ws.Range("N136").Value = -15696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10182.728
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H132").Value = 2380.3044
$ws.Range("I132").Value = 2156.6924
$ws.Range("J132").Value = 2671
$ws.Range("K132").Value = 6470.0772
$ws.Range("L132").Value = 8013
$ws.Range("M132").Value = -3940.0772
$ws.Range("N132").Value = -13073

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6221.2104
$ws.Range("I136").Value = 9266.583000000001
$ws.Range("J136").Value = 1000.5714
$ws.Range("K136").Value = 27799.749
$ws.Range("L136").Value = 3001.7142
$ws.Range("M136").Value = -25249.749
$ws.Range("N136").Value = -8101.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 807.48
$ws.Range("I96").Value = 690.1111
$ws.Range("K96").Value = 690.1111
$ws.Range("M96").Value = 682.8889
$ws.Range("H107").Value = 371.63635
$ws.Range("I107").Value = 398.66666
$ws.Range("K107").Value = 1195.99998
$ws.Range("M107").Value = 724.0000199999999
$ws.Range("H132").Value = 5032.8887
$ws.Range("I132").Value = 4216.5
$ws.Range("J132").Value = 6665.6665
$ws.Range("K132").Value = 12649.5
$ws.Range("L132").Value = 19996.9995
$ws.Range("M132").Value = -10119.5
$ws.Range("N132").Value = -25056.9995
